$p = $ppt.ActivePresentation

# --- Append 9 new blank slides at the end of the deck (final positions 11-19) ---
# ppLayoutBlank = 12
for ($i = 0; $i -lt 9; $i++) {
    $idx = $p.Slides.Count + 1
    $null = $p.Slides.Add($idx, 12)
}

# --- Draw the little "faces and an owl" doodle on the 7th new slide (slide 17) ---
$doodle = $p.Slides.Item(17)

function Style-Doodle($shp) {
    $shp.TextFrame.TextRange.Text = ""
    $shp.TextFrame.TextRange.ParagraphFormat.Alignment = 2   # ppAlignCenter
    $shp.TextFrame.VerticalAnchor = 3                        # msoAnchorMiddle
}

# msoShapeRectangle = 1
$rect1 = $doodle.Shapes.AddShape(1, 154.01748031496064, 146.4432283464567, 85.0488188976378, 11.339842519685039)
$rect1.Name = "Rectangle 1"
$rect1.Rotation = 26.198416666666667
Style-Doodle $rect1

$rect2 = $doodle.Shapes.AddShape(1, 260.13023622047245, 145.63094488188977, 96.38866141732284, 11.339842519685039)
$rect2.Name = "Rectangle 2"
$rect2.Rotation = 330.42401666666666
Style-Doodle $rect2

# msoShapeOval = 9
$oval1 = $doodle.Shapes.AddShape(9, 150.21291338582677, 173.61133858267718, 62.369133858267716, 39.68944881889764)
$oval1.Name = "Oval 3"
Style-Doodle $oval1

$oval2 = $doodle.Shapes.AddShape(9, 297.63086614173227, 173.61133858267718, 73.70897637795275, 45.359370078740156)
$oval2.Name = "Oval 4"
Style-Doodle $oval2

# msoShapeRoundedRectangle = 5
$round1 = $doodle.Shapes.AddShape(5, 218.251968503937, 298.3496062992126, 73.70897637795275, 85.0488188976378)
$round1.Name = "Rounded Rectangle 5"
Style-Doodle $round1
